# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.642.09"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3
$ws.Range("D3").Value = "3.492.71"
$ws.Range("E3").Value = "  +1.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.18"
$ws.Range("E5").Value = "  +1.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.29"
$ws.Range("E6").Value = "  +2.98%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  +11.24%  "

# Row 9
$ws.Range("D9").Value = "3.496.93"
$ws.Range("E9").Value = "  +1.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("E10").Value = "  -1.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +2.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  +2.05%  "

# Row 13
$ws.Range("D13").Value = "4.094.45"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14
$ws.Range("E14").Value = "  +0.23%  "

# Row 15
$ws.Range("E15").Value = "  +1.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.72"
$ws.Range("E16").Value = "  +5.62%  "

# Row 17
$ws.Range("D17").Value = "65.626.21"
$ws.Range("E17").Value = "  +2.73%  "

# Row 18
$ws.Range("D18").Value = "3.496.31"
$ws.Range("E18").Value = "  +3.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  +2.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.34"
$ws.Range("E20").Value = "  +0.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.00"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.27"
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("E23").Value = "  +3.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.37"
$ws.Range("E24").Value = "  +1.44%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.33%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("E26").Value = "  +4.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("E27").Value = "  +7.44%  "

# Row 28
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  +4.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +6.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  +2.87%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.74"
$ws.Range("E33").Value = "  +1.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.22"
$ws.Range("E34").Value = "  +4.77%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +7.90%  "

# Row 36
$ws.Range("E36").Value = "  +2.15%  "

# Row 37
$ws.Range("E37").Value = "  +6.21%  "

# Row 38
$ws.Range("D38").Value = "3.057.57"
$ws.Range("E38").Value = "  +5.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0777"
$ws.Range("E39").Value = "  +0.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.34"
$ws.Range("E40").Value = "  +1.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0325"
$ws.Range("E41").Value = "  +2.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.58"
$ws.Range("E42").Value = "  +3.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.16"
$ws.Range("E43").Value = "  +4.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.54"
$ws.Range("E44").Value = "  +1.11%  "

# Row 45
$ws.Range("E45").Value = "  +1.65%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.99"
$ws.Range("E46").Value = "  +9.71%  "

# Row 47
$ws.Range("E47").Value = "  +3.31%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "318.40"
$ws.Range("E48").Value = "  +8.53%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.23"
$ws.Range("E49").Value = "  +0.80%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.75"
$ws.Range("E50").Value = "  +4.36%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.109"
$ws.Range("E51").Value = "  +5.65%  "
